# Updates cryptos list prices and 1h volume percentages (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (plain decimals like "299.22") are forced to stay text: mark the cell as
# Text-formatted, assign the literal string, then clear the format again so
# the cell keeps its original (default) style, matching the source data which
# stores every Price/Volume value as inline text.

$ws.Range("D2").Value = "45.760.52"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.520.27"
$ws.Range("E3").Value = "  +10.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.22"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.32"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.588"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.34%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.95"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0803"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.67"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +9.63%  "
$ws.Range("D13").Value = "2.909.70"
$ws.Range("E13").Value = "  +11.02%  "
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").Value = "2.516.94"
$ws.Range("E15").Value = "  +9.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.878"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +10.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.55"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +7.38%  "
$ws.Range("D18").Value = "45.903.70"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.30"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +6.85%  "
$ws.Range("D20").Value = "0.0₃0969"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.50"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +11.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.78"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.76"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.87"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +9.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.90"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.71%  "
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.68"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +14.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +16.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.79"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.92%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.79"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.19"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +31.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.18"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0805"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.04%  "
$ws.Range("E37").Value = "  +5.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.118"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.80"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.11"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0312"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.38"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.20%  "
$ws.Range("D43").Value = "2.017.95"
$ws.Range("E43").Value = "  +11.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.31"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.65"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +24.26%  "
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.94"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +12.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.79"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +13.76%  "
$ws.Range("D50").Value = "2.773.17"
$ws.Range("E50").Value = "  +10.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.194"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.63%  "
